$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.651.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.407.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.404.28"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.118"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -11.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.99"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -11.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -8.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.987.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000176"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -10.87%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.115"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.02%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.10"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.26%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.403.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.674.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -14.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "380.40"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.545"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -9.10%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.97"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.547.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -11.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.424.09"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.42%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.88"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.77%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.140"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "167.81"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.67"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -11.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -12.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.59"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -12.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0746"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.30%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -15.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.59"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -10.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.40"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.98%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.148.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.85%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -15.56%  "
